$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dimension / header row: add week 51 (BB1) and week 52 (BC1)
$ws.Range("BB1").Value = "51"
$ws.Range("BB1").Font.Bold = $true
$ws.Range("BB1").HorizontalAlignment = -4108
$ws.Range("BC1").Value = "52"
$ws.Range("BC1").Font.Bold = $true
$ws.Range("BC1").HorizontalAlignment = -4108

# Weekly data for week 51 (BB) and week 52 (BC) per row
$ws.Range("BB2").Value = 0
$ws.Range("BC2").Value = 0
$ws.Range("BB3").Value = 0
$ws.Range("BC3").Value = 0
$ws.Range("BB5").Value = 0
$ws.Range("BC5").Value = 0
$ws.Range("BB6").Value = 25
$ws.Range("BC6").Value = 13
$ws.Range("BB7").Value = 13
$ws.Range("BC7").Value = 28
$ws.Range("BB8").Value = 6
$ws.Range("BC8").Value = 11
$ws.Range("BB9").Value = 0
$ws.Range("BC9").Value = 0
$ws.Range("BB10").Value = 0
$ws.Range("BC10").Value = 0
$ws.Range("BB11").Value = 0
$ws.Range("BB14").Value = 0
$ws.Range("BC14").Value = 0
$ws.Range("BB15").Value = 0
$ws.Range("BB16").Value = 0
$ws.Range("BC16").Value = 0
$ws.Range("BB17").Value = 0
$ws.Range("BB23").Value = 0
$ws.Range("BB25").Value = 2
$ws.Range("BC25").Value = 0
$ws.Range("BB28").Value = 4
$ws.Range("BC28").Value = 3
$ws.Range("BB29").Value = 0
$ws.Range("BC29").Value = 0
$ws.Range("BB30").Value = 8
$ws.Range("BC30").Value = 13
$ws.Range("BB31").Value = 0
$ws.Range("BC31").Value = 0
$ws.Range("BB35").Value = 2
$ws.Range("BC35").Value = 4
$ws.Range("BB36").Value = 0
$ws.Range("BC36").Value = 0
$ws.Range("BB37").Value = 0
$ws.Range("BC37").Value = 0
$ws.Range("BB38").Value = 0
$ws.Range("BC38").Value = 0
$ws.Range("BB41").Value = 0
$ws.Range("BC41").Value = 0
$ws.Range("BB42").Value = 0
$ws.Range("BC42").Value = 0
$ws.Range("BB43").Value = 0
$ws.Range("BC43").Value = 0
$ws.Range("BB45").Value = 0
$ws.Range("BC45").Value = 0
$ws.Range("BB46").Value = 0
$ws.Range("BC46").Value = 0
$ws.Range("BB47").Value = 0
$ws.Range("BC47").Value = 0
$ws.Range("BB48").Value = 0
$ws.Range("BC48").Value = 0
$ws.Range("BB49").Value = 0
$ws.Range("BC49").Value = 0
$ws.Range("BB50").Value = 0
$ws.Range("BC50").Value = 0
$ws.Range("BB51").Value = 0
$ws.Range("BC51").Value = 0
$ws.Range("BB54").Value = 0
$ws.Range("BC54").Value = 0
$ws.Range("BB55").Value = 0
$ws.Range("BC55").Value = 0
$ws.Range("BB56").Value = 0
$ws.Range("BC56").Value = 0
$ws.Range("BB57").Value = 0
$ws.Range("BC57").Value = 0
$ws.Range("BB58").Value = 0
$ws.Range("BC58").Value = 0
$ws.Range("BB59").Value = 0
$ws.Range("BC59").Value = 0
